# ExamPortal sample workbook: replace the 4-question sample sheet with a
# single "richest man in the world" single-option sample question, and
# drop the old hyperlink sample cell's link (keeping only its formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Remove the old hyperlink on H4 ----
$ws.Hyperlinks.Delete()

# ---- Wipe out rows 2-4 (old sample questions) ----
$ws.Range("A2:I4").Clear()

# ---- New single sample question on row 1 ----
$ws.Range("A1").Value = "who is the richest man in the world?"
$ws.Range("B1").Value = "Mukesh Ambani"
$ws.Range("C1").Value = "Bill Gates"
$ws.Range("D1").Value = "Jeff Bezos"
$ws.Range("E1").Value = "Warren Buffet"
$ws.Range("F1").Value = " option3"
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = "null"
$ws.Range("I1").Value = "singleOption"

# ---- Re-apply the (now link-less) hyperlink look to the sample cell ----
$ws.Range("H4").Style = "Hyperlink"

# ---- New column widths for the answer-option columns ----
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 17.59
$ws.Range("D1").EntireColumn.ColumnWidth = 10.75
$ws.Range("E1").EntireColumn.ColumnWidth = 15.25

# ---- Updated view/selection state ----
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I2").Select()
